$wb = $excel.ActiveWorkbook

# --- "About" sheet: update the revision/last-updated date (C1) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- "RAF-capacity" sheet: hydrogen rows get full capacity credit (1 instead of 0.3) ---
$wsCap = $wb.Worksheets.Item("RAF-capacity")
$wsCap.Range("B24").Value = 1
$wsCap.Range("B25").Value = 1

# Widen column A slightly on this sheet (closest achievable quantized width to the
# recorded 29.04296875 target given this engine's column-width rounding)
$wsCap.Columns.Item(1).ColumnWidth = 28.14

# --- View state: the workbook was left with RAF-capacity as the active/visible tab,
# scrolled down so row 14 is at the top, zoomed to 80%, with B25 selected ---
$wsCap.Activate()
$excel.ActiveWindow.Zoom = 80
$wsCap.Range("B25").Select()
